$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1400
$ws.Range("I69").Value = 1400
$ws.Range("K69").Value = 4200
$ws.Range("M69").Value = -3326

$ws.Range("H72").Value = 1400
$ws.Range("I72").Value = 1400
$ws.Range("K72").Value = 12600
$ws.Range("M72").Value = -8232

$ws.Range("H86").Value = 11355
$ws.Range("I86").Value = 17441.666
$ws.Range("J86").Value = 2225
$ws.Range("K86").Value = 17441.666
$ws.Range("L86").Value = 2225
$ws.Range("M86").Value = -16318.666
$ws.Range("N86").Value = -4471

$ws.Range("H89").Value = 11355
$ws.Range("I89").Value = 17441.666
$ws.Range("J89").Value = 2225
$ws.Range("K89").Value = 87208.33
$ws.Range("L89").Value = 11125
$ws.Range("M89").Value = -81592.33
$ws.Range("N89").Value = -22357

$ws.Range("H129").Value = 1125.0834
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1125.0834
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3375.2502
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -13375.2502

$ws.Range("H137").Value = 1175.6428
$ws.Range("I137").Value = 767.7931
$ws.Range("J137").Value = 2085.4614
$ws.Range("K137").Value = 2303.3793
$ws.Range("L137").Value = 6256.3842
$ws.Range("M137").Value = 246.6206999999999
$ws.Range("N137").Value = -11356.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9770.272000000001
$ws.Range("I32").Value = 9645.621999999999
$ws.Range("J32").Value = 10429.143
$ws.Range("K32").Value = 9645.621999999999
$ws.Range("L32").Value = 10429.143
$ws.Range("M32").Value = -9358.621999999999
$ws.Range("N32").Value = -11003.143

$ws.Range("H132").Value = 3726.2683
$ws.Range("I132").Value = 987.6896400000001
$ws.Range("J132").Value = 10344.5
$ws.Range("K132").Value = 2963.06892
$ws.Range("L132").Value = 31033.5
$ws.Range("M132").Value = -433.0689200000002
$ws.Range("N132").Value = -36093.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 4010.3823
$ws.Range("I80").Value = 843.4286
$ws.Range("J80").Value = 6227.25
$ws.Range("K80").Value = 843.4286
$ws.Range("L80").Value = 6227.25
$ws.Range("M80").Value = 154.5714
$ws.Range("N80").Value = -8223.25

$ws.Range("H83").Value = 4010.3823
$ws.Range("I83").Value = 843.4286
$ws.Range("J83").Value = 6227.25
$ws.Range("K83").Value = 4217.143
$ws.Range("L83").Value = 31136.25
$ws.Range("M83").Value = 774.857
$ws.Range("N83").Value = -41120.25

$ws.Range("H105").Value = 34485372
$ws.Range("I105").Value = 2757.4092
$ws.Range("K105").Value = 2757.4092
$ws.Range("M105").Value = -1010.4092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 23590
$ws.Range("I56").Value = 8100
$ws.Range("K56").Value = 8100
$ws.Range("M56").Value = -7255

$ws.Range("H58").Value = 18868550
$ws.Range("I58").Value = 28571948
$ws.Range("J58").Value = 832.2222
$ws.Range("K58").Value = 28571948
$ws.Range("L58").Value = 832.2222
$ws.Range("M58").Value = -28571745
$ws.Range("N58").Value = -1238.2222

$ws.Range("H132").Value = 12825835
$ws.Range("I132").Value = 1139.9412
$ws.Range("J132").Value = 37050260
$ws.Range("K132").Value = 3419.8236
$ws.Range("L132").Value = 111150780
$ws.Range("M132").Value = -889.8235999999997
$ws.Range("N132").Value = -111155840

$ws.Range("H136").Value = 18868550
$ws.Range("I136").Value = 28571948
$ws.Range("J136").Value = 832.2222
$ws.Range("K136").Value = 85715844
$ws.Range("L136").Value = 2496.6666
$ws.Range("M136").Value = -85713294
$ws.Range("N136").Value = -7596.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7144804.5
$ws.Range("I5").Value = 283.82977
$ws.Range("J5").Value = 21744478
$ws.Range("K5").Value = 851.4893099999999
$ws.Range("L5").Value = 65233434
$ws.Range("M5").Value = -739.4893099999999
$ws.Range("N5").Value = -65233658

$ws.Range("H107").Value = 52637876
$ws.Range("I107").Value = 360
$ws.Range("K107").Value = 1080
$ws.Range("M107").Value = 840

$ws.Range("H109").Value = 3109.0952
$ws.Range("I109").Value = 616.44446
$ws.Range("J109").Value = 4978.5835
$ws.Range("K109").Value = 1849.33338
$ws.Range("L109").Value = 14935.7505
$ws.Range("M109").Value = -809.33338
$ws.Range("N109").Value = -17015.7505

$ws.Range("H122").Value = 15628722
$ws.Range("I122").Value = 27778080
$ws.Range("J122").Value = 8119.9287
$ws.Range("K122").Value = 250002720
$ws.Range("L122").Value = 73079.35830000001
$ws.Range("M122").Value = -250000270
$ws.Range("N122").Value = -77979.35830000001

$ws.Range("H131").Value = 709.24
$ws.Range("I131").Value = 401.58334
$ws.Range("J131").Value = 751.1932
$ws.Range("K131").Value = 1204.75002
$ws.Range("L131").Value = 2253.5796
$ws.Range("M131").Value = 3835.24998
$ws.Range("N131").Value = -12333.5796

$ws.Range("H132").Value = 14289962
$ws.Range("J132").Value = 19236146
$ws.Range("L132").Value = 173125314
$ws.Range("N132").Value = -173130374

$ws.Range("H135").Value = 7144804.5
$ws.Range("I135").Value = 283.82977
$ws.Range("J135").Value = 21744478
$ws.Range("K135").Value = 2554.46793
$ws.Range("L135").Value = 195700302
$ws.Range("M135").Value = -19.4679299999998
$ws.Range("N135").Value = -195705372

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 19700
$ws.Range("J68").Value = 19700
$ws.Range("L68").Value = 19700
$ws.Range("N68").Value = -21322

$ws.Range("H71").Value = 19700
$ws.Range("J71").Value = 19700
$ws.Range("L71").Value = 59100
$ws.Range("N71").Value = -67212

$ws.Range("H80").Value = 12504301
$ws.Range("I80").Value = 5349.8335
$ws.Range("J80").Value = 50001150
$ws.Range("K80").Value = 5349.8335
$ws.Range("L80").Value = 50001150
$ws.Range("M80").Value = -4351.8335
$ws.Range("N80").Value = -50003146

$ws.Range("H83").Value = 12504301
$ws.Range("I83").Value = 5349.8335
$ws.Range("J83").Value = 50001150
$ws.Range("K83").Value = 26749.1675
$ws.Range("L83").Value = 250005750
$ws.Range("M83").Value = -21757.1675
$ws.Range("N83").Value = -250015734

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7794.9644
$ws.Range("I132").Value = 2550.5
$ws.Range("J132").Value = 11728.3125
$ws.Range("K132").Value = 7651.5
$ws.Range("L132").Value = 35184.9375
$ws.Range("M132").Value = -5121.5
$ws.Range("N132").Value = -40244.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2419.2466
$ws.Range("I136").Value = 2762.5
$ws.Range("J136").Value = 1673.0435
$ws.Range("K136").Value = 8287.5
$ws.Range("L136").Value = 5019.1305
$ws.Range("M136").Value = -5737.5
$ws.Range("N136").Value = -10119.1305
